$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 343.2857
$ws.Range("I2").Value = 357
$ws.Range("K2").Value = 357
$ws.Range("M2").Value = -244
$ws.Range("H6").Value = 303.30768
$ws.Range("I6").Value = 138.18182
$ws.Range("K6").Value = 414.5454599999999
$ws.Range("M6").Value = -302.5454599999999
$ws.Range("H9").Value = 167.1875
$ws.Range("I9").Value = 155.42857
$ws.Range("J9").Value = 249.5
$ws.Range("K9").Value = 155.42857
$ws.Range("L9").Value = 249.5
$ws.Range("M9").Value = 13.57142999999999
$ws.Range("N9").Value = -587.5
$ws.Range("H33").Value = 779.55554
$ws.Range("I33").Value = 294.23077
$ws.Range("K33").Value = 294.23077
$ws.Range("M33").Value = -65.23077
$ws.Range("H38").Value = 3158.2727
$ws.Range("I38").Value = 2405.125
$ws.Range("K38").Value = 7215.375
$ws.Range("M38").Value = -6843.375
$ws.Range("H43").Value = 7374.5
$ws.Range("I43").Value = 3500
$ws.Range("K43").Value = 3500
$ws.Range("M43").Value = -3431
$ws.Range("H47").Value = 6420.5713
$ws.Range("I47").Value = 5832.3335
$ws.Range("J47").Value = 9950
$ws.Range("K47").Value = 5832.3335
$ws.Range("L47").Value = 9950
$ws.Range("M47").Value = -4860.3335
$ws.Range("N47").Value = -11894
$ws.Range("H51").Value = 9397.571
$ws.Range("I51").Value = 9397
$ws.Range("K51").Value = 9397
$ws.Range("M51").Value = -8913
$ws.Range("H58").Value = 1616.7059
$ws.Range("I58").Value = 407.63635
$ws.Range("J58").Value = 3833.3333
$ws.Range("K58").Value = 1222.90905
$ws.Range("L58").Value = 11499.9999
$ws.Range("M58").Value = -1072.90905
$ws.Range("N58").Value = -11799.9999
$ws.Range("H112").Value = 73356.21
$ws.Range("I112").Value = 866.3333
$ws.Range("K112").Value = 2598.9999
$ws.Range("M112").Value = -1490.9999
$ws.Range("H116").Value = 12527921
$ws.Range("I116").Value = 25048242
$ws.Range("K116").Value = 25048242
$ws.Range("M116").Value = -25044800
$ws.Range("H130").Value = 58853
$ws.Range("J130").Value = 58853
$ws.Range("L130").Value = 58853
$ws.Range("N130").Value = -68893
$ws.Range("H132").Value = 3120.262
$ws.Range("I132").Value = 2969.639
$ws.Range("J132").Value = 4024
$ws.Range("K132").Value = 8908.917000000001
$ws.Range("L132").Value = 12072
$ws.Range("M132").Value = -6378.917000000001
$ws.Range("N132").Value = -17132
$ws.Range("H138").Value = 200122.67
$ws.Range("J138").Value = 4692.913
$ws.Range("L138").Value = 14078.739
$ws.Range("N138").Value = -24358.739

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H9").Value = 20000
$ws.Range("J9").Value = 20000
$ws.Range("L9").Value = 20000
$ws.Range("N9").Value = -20340
$ws.Range("H20").Value = 20000
$ws.Range("J20").Value = 20000
$ws.Range("L20").Value = 20000
$ws.Range("N20").Value = -20540
$ws.Range("H32").Value = 2907.203
$ws.Range("I32").Value = 2470.754
$ws.Range("K32").Value = 2470.754
$ws.Range("M32").Value = -2183.754
$ws.Range("H61").Value = 6486.7036
$ws.Range("J61").Value = 5082
$ws.Range("L61").Value = 5082
$ws.Range("N61").Value = -5506
$ws.Range("H63").Value = 5234.3335
$ws.Range("I63").Value = 3872.7144
$ws.Range("K63").Value = 3872.7144
$ws.Range("M63").Value = -3186.7144
$ws.Range("H66").Value = 5234.3335
$ws.Range("I66").Value = 3872.7144
$ws.Range("K66").Value = 19363.572
$ws.Range("M66").Value = -15931.572
$ws.Range("H74").Value = 3515.1177
$ws.Range("I74").Value = 983.8
$ws.Range("J74").Value = 22500
$ws.Range("K74").Value = 983.8
$ws.Range("L74").Value = 22500
$ws.Range("M74").Value = -109.8
$ws.Range("N74").Value = -24248
$ws.Range("H77").Value = 3515.1177
$ws.Range("I77").Value = 983.8
$ws.Range("J77").Value = 22500
$ws.Range("K77").Value = 4919
$ws.Range("L77").Value = 112500
$ws.Range("M77").Value = -551
$ws.Range("N77").Value = -121236
$ws.Range("H97").Value = 17449.762
$ws.Range("I97").Value = 7090.294
$ws.Range("K97").Value = 7090.294
$ws.Range("M97").Value = -6594.294
$ws.Range("H136").Value = 6486.7036
$ws.Range("J136").Value = 5082
$ws.Range("L136").Value = 15246
$ws.Range("N136").Value = -20346

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 7845.643
$ws.Range("J64").Value = 3666.889
$ws.Range("L64").Value = 3666.889
$ws.Range("N64").Value = -4116.889
$ws.Range("H67").Value = 7845.643
$ws.Range("J67").Value = 3666.889
$ws.Range("L67").Value = 3666.889
$ws.Range("N67").Value = -5226.889

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 6054.4
$ws.Range("J86").Value = 8712.333
$ws.Range("L86").Value = 8712.333
$ws.Range("N86").Value = -10958.333
$ws.Range("H89").Value = 6054.4
$ws.Range("J89").Value = 8712.333
$ws.Range("L89").Value = 43561.665
$ws.Range("N89").Value = -54793.665
$ws.Range("H92").Value = 74994.5
$ws.Range("J92").Value = 74994.5
$ws.Range("L92").Value = 74994.5
$ws.Range("N92").Value = -79986.5
$ws.Range("H109").Value = 83884
$ws.Range("J109").Value = 83884
$ws.Range("L109").Value = 83884
$ws.Range("N109").Value = -85964

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 250000130
$ws.Range("I7").Value = 150
$ws.Range("J7").Value = 500000100
$ws.Range("K7").Value = 450
$ws.Range("L7").Value = 1500000300
$ws.Range("M7").Value = -338
$ws.Range("N7").Value = -1500000524
$ws.Range("H120").Value = 66674668
$ws.Range("I120").Value = 111117780
$ws.Range("K120").Value = 333353340
$ws.Range("M120").Value = -333348502

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H47").Value = 17000
$ws.Range("I47").Value = 9000
$ws.Range("K47").Value = 9000
$ws.Range("M47").Value = -8432
$ws.Range("H80").Value = 11631.272
$ws.Range("I80").Value = 14312.5
$ws.Range("J80").Value = 4481.3335
$ws.Range("K80").Value = 14312.5
$ws.Range("L80").Value = 4481.3335
$ws.Range("M80").Value = -13314.5
$ws.Range("N80").Value = -6477.3335
$ws.Range("H83").Value = 11631.272
$ws.Range("I83").Value = 14312.5
$ws.Range("J83").Value = 4481.3335
$ws.Range("K83").Value = 71562.5
$ws.Range("L83").Value = 22406.6675
$ws.Range("M83").Value = -66570.5
$ws.Range("N83").Value = -32390.6675
$ws.Range("H113").Value = 21449.5
$ws.Range("I113").Value = 21449.5
$ws.Range("K113").Value = 21449.5
$ws.Range("M113").Value = -19279.5
$ws.Range("H132").Value = 3059.6584
$ws.Range("I132").Value = 2491.8
$ws.Range("J132").Value = 6372.1665
$ws.Range("K132").Value = 7475.400000000001
$ws.Range("L132").Value = 19116.4995
$ws.Range("M132").Value = -4945.400000000001
$ws.Range("N132").Value = -24176.4995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2868.9546
$ws.Range("I16").Value = 2255.85
$ws.Range("K16").Value = 2255.85
$ws.Range("M16").Value = -2085.85
$ws.Range("H68").Value = 6057.364
$ws.Range("J68").Value = 5300
$ws.Range("L68").Value = 5300
$ws.Range("N68").Value = -6798
$ws.Range("H71").Value = 6057.364
$ws.Range("J71").Value = 5300
$ws.Range("L71").Value = 26500
$ws.Range("N71").Value = -33988
$ws.Range("H136").Value = 4143.8184
$ws.Range("I136").Value = 1366.68
$ws.Range("K136").Value = 4100.04
$ws.Range("M136").Value = -1550.04

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 12065.4
$ws.Range("I81").Value = 16178.1
$ws.Range("J81").Value = 3840
$ws.Range("K81").Value = 32356.2
$ws.Range("L81").Value = 7680
$ws.Range("M81").Value = -31295.2
$ws.Range("N81").Value = -9802
$ws.Range("H84").Value = 12065.4
$ws.Range("I84").Value = 16178.1
$ws.Range("J84").Value = 3840
$ws.Range("K84").Value = 161781
$ws.Range("L84").Value = 38400
$ws.Range("M84").Value = -156477
$ws.Range("N84").Value = -49008
$ws.Range("H122").Value = 5162.8335
$ws.Range("I122").Value = 3494.7273
$ws.Range("K122").Value = 10484.1819
$ws.Range("M122").Value = -8034.1819
$ws.Range("H136").Value = 2163.5938
$ws.Range("I136").Value = 1457.64
$ws.Range("J136").Value = 4684.857
$ws.Range("K136").Value = 4372.92
$ws.Range("L136").Value = 14054.571
$ws.Range("M136").Value = -1822.92
$ws.Range("N136").Value = -19154.571
